$d = $word.ActiveDocument

# Update the date
$d.Content.Find.Execute("2023-07-21", $true, $false, $false, $false, $false,
                         $true, 1, $false, "8/9/23", 2)

# Replace all occurrences of HWIs with GIWs
$d.Content.Find.Execute("HWIs", $true, $false, $false, $false, $false,
                         $true, 1, $false, "GIWs", 2)
